$d = $word.ActiveDocument

# Replace first paragraph's text
$d.Paragraphs(1).Range.Text = "Observation of a new particle in the search for the Standard Model Higgs boson with the ATLAS detector at the LHC`r"

# Replace second paragraph's text (no leading space in target)
$d.Paragraphs(2).Range.Text = "PYTHIA 6.4 physics and manual`r"

# Delete paragraphs 3 through 8 (trailing paragraphs no longer present)
for ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {
    $d.Paragraphs($i).Range.Delete()
}
